# Updated cryptos list on Tue Jul 16 19:54:29 UTC 2024 with GitHub Actions
# Refresh Price (column D) / Volume(1h) (column E) figures, and restore the
# correct row order for two coin pairs (OKB/VeChain and dogwifhat/Stellar)
# whose rank (column A) stayed fixed while their B/C/D/E data swapped.
#
# Note: several new Price values look numeric (e.g. "577.73"); a leading
# apostrophe forces Excel to store them as text, matching the original
# inline-string/text representation used throughout column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.895.10'
$ws.Range('E2').Value = '  +2.42%  '

$ws.Range('D3').Value = '3.466.26'
$ws.Range('E3').Value = '  +1.98%  '

$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').Value = '''577.73'
$ws.Range('E5').Value = '  +0.22%  '

$ws.Range('D6').Value = '''161.78'
$ws.Range('E6').Value = '  +4.40%  '

$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('D8').Value = '3.467.36'
$ws.Range('E8').Value = '  +1.78%  '

$ws.Range('E9').Value = '  +8.69%  '

$ws.Range('E10').Value = '  -2.43%  '

$ws.Range('E11').Value = '  +3.77%  '

$ws.Range('D12').Value = '''0.441'
$ws.Range('E12').Value = '  +1.65%  '

$ws.Range('D13').Value = '4.063.42'
$ws.Range('E13').Value = '  +1.91%  '

$ws.Range('E14').Value = '  -2.83%  '

$ws.Range('E15').Value = '  +5.48%  '

$ws.Range('D16').Value = '''28.18'
$ws.Range('E16').Value = '  +4.16%  '

$ws.Range('D17').Value = '64.883.10'
$ws.Range('E17').Value = '  +2.15%  '

$ws.Range('D18').Value = '3.455.72'
$ws.Range('E18').Value = '  +2.52%  '

$ws.Range('E19').Value = '  +0.12%  '

$ws.Range('D20').Value = '''14.34'
$ws.Range('E20').Value = '  +1.88%  '

$ws.Range('D21').Value = '''390.29'
$ws.Range('E21').Value = '  +0.63%  '

$ws.Range('D22').Value = '''8.19'
$ws.Range('E22').Value = '  -2.48%  '

$ws.Range('D23').Value = '''0.547'
$ws.Range('E23').Value = '  +2.24%  '

$ws.Range('D24').Value = '''72.86'
$ws.Range('E24').Value = '  +2.54%  '

$ws.Range('E25').Value = '  +0.20%  '

$ws.Range('E26').Value = '  +16.68%  '

$ws.Range('D27').Value = '''9.56'
$ws.Range('E27').Value = '  +0.37%  '

$ws.Range('D28').Value = '''0.181'
$ws.Range('E28').Value = '  -0.52%  '

$ws.Range('E29').Value = '  +0.17%  '

$ws.Range('D30').Value = '''6.22'
$ws.Range('E30').Value = '  +8.91%  '

$ws.Range('E31').Value = '  +7.16%  '

$ws.Range('D32').Value = '''2.05'
$ws.Range('E32').Value = '  +0.52%  '

$ws.Range('D33').Value = '''23.67'
$ws.Range('E33').Value = '  +1.90%  '

$ws.Range('D34').Value = '''6.53'
$ws.Range('E34').Value = '  -0.37%  '

$ws.Range('E35').Value = '  +0.15%  '

$ws.Range('D36').Value = '''7.08'
$ws.Range('E36').Value = '  +5.42%  '

$ws.Range('E37').Value = '  +1.22%  '

$ws.Range('D38').Value = '''162.02'
$ws.Range('E38').Value = '  +2.45%  '

$ws.Range('D39').Value = '''1.91'
$ws.Range('E39').Value = '  +1.85%  '

$ws.Range('D40').Value = '3.021.84'
$ws.Range('E40').Value = '  +3.79%  '

$ws.Range('E41').Value = '  -0.37%  '

$ws.Range('D42').Value = '''27.36'
$ws.Range('E42').Value = '  -0.96%  '

$ws.Range('E43').Value = '  +5.77%  '

$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').Value = '''0.0316'
$ws.Range('E44').Value = '  -0.78%  '

$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '''42.69'
$ws.Range('E45').Value = '  +3.04%  '

$ws.Range('E46').Value = '  +2.01%  '

$ws.Range('D47').Value = '''24.31'
$ws.Range('E47').Value = '  +8.29%  '

$ws.Range('E48').Value = '  +1.95%  '

$ws.Range('E49').Value = '  +8.11%  '

$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '''0.107'
$ws.Range('E50').Value = '  +4.02%  '

$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '''2.17'
$ws.Range('E51').Value = '  +9.85%  '
